$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the taxon-identifying data between row 4 and row 5
# (columns A, B, D, E, F, G, H, Q, R). All other columns in these two
# rows were already identical, so only these columns need to change.

$cols = 1, 2, 4, 5, 6, 7, 8, 17, 18

foreach ($col in $cols) {
    $cell4 = $ws.Cells.Item(4, $col)
    $cell5 = $ws.Cells.Item(5, $col)

    $v4 = $cell4.Value2
    $v5 = $cell5.Value2

    $cell4.Value2 = $v5
    $cell5.Value2 = $v4
}
